$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "square"
$ws.Range("B1").Value = "loc1"
$ws.Range("C1").Value = "loc2"
$ws.Range("D1").Value = "corrAns"

$ws.Range("D2").Select()
